$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 189, shifting the existing rows 189:303 down to 190:304.
$ws.Rows("189:189").Insert()

# Populate the newly inserted row 189 with the new weekly record.
$ws.Range("A189").Value = 4
$ws.Range("B189").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C189").Value = "Los Lagos"
$ws.Range("D189").Value = 44873
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 100112044
$ws.Range("G189").Value = "Perejil"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 180
$ws.Range("K189").Value = 5000
$ws.Range("L189").Value = 5000
$ws.Range("M189").Value = 5000
$ws.Range("N189").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O189").Value = "Región Metropolitana"
$ws.Range("P189").Value = 1667
$ws.Range("Q189").Value = 3
$ws.Range("R189").Value = "Hortaliza"
